# Update countries & provincias Spain
#
# Refreshes the COVID-19 "paises" data table (sheet "Pais") with a newer
# snapshot of figures (columns B:H = Casos totales, Nuevos casos, Casos
# activos, Recuperados, Casos criticos, Muertes hoy, Muertes). The table is
# kept sorted descending by column B ("Casos totales"), so several countries
# whose totals overtook their neighbour swap row position as part of this
# refresh (Panama/Belgica, Surinam/Republica de Africa Central/Ruanda,
# Guyana/Sudan del Sur, Montserrat/Islas Malvinas). The "Datos actualizados"
# timestamp banner in A1 is also bumped to the new refresh time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Banner timestamp (A1)
$ws.Range("A1").Value2 = "Datos actualizados a 26 de Septiembre de 2020 a las 01:44"

# Estados Unidos (row 4)
$ws.Range("B4").Value2 = 7233536
$ws.Range("C4").Value2 = 48314
$ws.Range("D4").Value2 = 4475387
$ws.Range("E4").Value2 = 2549868
$ws.Range("G4").Value2 = 761
$ws.Range("H4").Value2 = 208281

# Peru (row 9)
$ws.Range("B9").Value2 = 794584
$ws.Range("C9").Value2 = 5654
$ws.Range("D9").Value2 = 650948
$ws.Range("E9").Value2 = 111599
$ws.Range("G9").Value2 = 99
$ws.Range("H9").Value2 = 32037

# Argentina (row 12)
$ws.Range("B12").Value2 = 691235
$ws.Range("C12").Value2 = 12969
$ws.Range("E12").Value2 = 139438
$ws.Range("G12").Value2 = 442
$ws.Range("H12").Value2 = 15208

# Panama moves ahead of Belgica (rows 36-37 swap countries + data)
$ws.Range("A36").Value2 = "Panama"
$ws.Range("B36").Value2 = 109431
$ws.Range("C36").Value2 = 705
$ws.Range("D36").Value2 = 86158
$ws.Range("E36").Value2 = 20962
$ws.Range("G36").Value2 = 14
$ws.Range("H36").Value2 = 2311

$ws.Range("A37").Value2 = "Belgica"
$ws.Range("B37").Value2 = 108768
$ws.Range("C37").Value2 = 1881
$ws.Range("D37").Value2 = 19123
$ws.Range("E37").Value2 = 79680
$ws.Range("G37").Value2 = 6
$ws.Range("H37").Value2 = 9965

# Chequia (row 57)
$ws.Range("B57").Value2 = 61318
$ws.Range("C57").Value2 = 2944
$ws.Range("D57").Value2 = 30740
$ws.Range("E57").Value2 = 29997
$ws.Range("G57").Value2 = 14
$ws.Range("H57").Value2 = 581

# Nigeria (row 58)
$ws.Range("B58").Value2 = 58062
$ws.Range("C58").Value2 = 213
$ws.Range("D58").Value2 = 49606
$ws.Range("E58").Value2 = 7353
$ws.Range("G58").Value2 = 1
$ws.Range("H58").Value2 = 1103

# Bulgaria (row 84)
$ws.Range("B84").Value2 = 19828
$ws.Range("C84").Value2 = 255
$ws.Range("D84").Value2 = 14132
$ws.Range("E84").Value2 = 4907
$ws.Range("G84").Value2 = 4
$ws.Range("H84").Value2 = 789

# Noruega (row 93)
$ws.Range("B93").Value2 = 13545
$ws.Range("C93").Value2 = 139
$ws.Range("E93").Value2 = 2904

# Gabon (row 106)
$ws.Range("B106").Value2 = 8728
$ws.Range("C106").Value2 = 12
$ws.Range("D106").Value2 = 7934
$ws.Range("E106").Value2 = 740

# Surinam moves ahead of Republica de Africa Central and Ruanda
# (rows 126-128 take the data of the row that used to be one below them)
$ws.Range("A126").Value2 = "Surinam"
$ws.Range("B126").Value2 = 4817
$ws.Range("C126").Value2 = 28
$ws.Range("D126").Value2 = 4596
$ws.Range("E126").Value2 = 119
$ws.Range("H126").Value2 = 102

$ws.Range("B127").Value2 = 4806
$ws.Range("C127").Value2 = 2
$ws.Range("D127").Value2 = 1840
$ws.Range("E127").Value2 = 2904
$ws.Range("G127").Value2 = 0
$ws.Range("H127").Value2 = 62

$ws.Range("A128").Value2 = "Republica de Africa Central"
$ws.Range("B128").Value2 = 4798
$ws.Range("C128").Value2 = 9
$ws.Range("D128").Value2 = 3080
$ws.Range("E128").Value2 = 1689
$ws.Range("G128").Value2 = 2
$ws.Range("H128").Value2 = 29

# Guyana moves ahead of Sudan del Sur (rows 147-148 swap countries + data)
$ws.Range("A147").Value2 = "Guyana"
$ws.Range("B147").Value2 = 2709
$ws.Range("C147").Value2 = 130
$ws.Range("D147").Value2 = 1490
$ws.Range("E147").Value2 = 1146
$ws.Range("G147").Value2 = 2
$ws.Range("H147").Value2 = 73

$ws.Range("A148").Value2 = "Sudan del Sur"
$ws.Range("B148").Value2 = 2676
$ws.Range("C148").Value2 = 7
$ws.Range("D148").Value2 = 1290
$ws.Range("E148").Value2 = 1337
$ws.Range("H148").Value2 = 49

# Uruguay (row 154)
$ws.Range("B154").Value2 = 1967
$ws.Range("C154").Value2 = 8
$ws.Range("D154").Value2 = 1710
$ws.Range("E154").Value2 = 210

# Santo Tome y Principe (row 169)
$ws.Range("B169").Value2 = 911
$ws.Range("C169").Value2 = 1
$ws.Range("E169").Value2 = 15

# Montserrat moves ahead of Islas Malvinas (rows 215-216 swap countries + data)
$ws.Range("A215").Value2 = "Montserrat"
$ws.Range("D215").Value2 = 12
$ws.Range("H215").Value2 = 1

$ws.Range("A216").Value2 = "Islas Malvinas"
$ws.Range("D216").Value2 = 13
$ws.Range("H216").Value2 = 0
